# Apply market-data value updates across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 179.6
$ws.Range("I12").Value = 166.33333
$ws.Range("K12").Value = 166.33333
$ws.Range("M12").Value = 3.666670000000011

$ws.Range("H13").Value = 7028.625
$ws.Range("J13").Value = 9825
$ws.Range("L13").Value = 9825
$ws.Range("N13").Value = -10163

$ws.Range("H42").Value = 66667584
$ws.Range("I42").Value = 90910264
$ws.Range("J42").Value = 214.25
$ws.Range("K42").Value = 272730792
$ws.Range("L42").Value = 642.75
$ws.Range("M42").Value = -272730562
$ws.Range("N42").Value = -1102.75

$ws.Range("H103").Value = 865.0833
$ws.Range("I103").Value = 696.25
$ws.Range("J103").Value = 949.5
$ws.Range("K103").Value = 2088.75
$ws.Range("L103").Value = 2848.5
$ws.Range("M103").Value = -1502.75
$ws.Range("N103").Value = -4020.5

$ws.Range("H116").Value = 9944.571
$ws.Range("I116").Value = 8674.857
$ws.Range("K116").Value = 8674.857
$ws.Range("M116").Value = -5232.857

$ws.Range("H123").Value = 63587.25
$ws.Range("J123").Value = 63587.25
$ws.Range("L123").Value = 63587.25
$ws.Range("N123").Value = -73387.25

$ws.Range("H131").Value = 11910892
$ws.Range("I131").Value = 23811384
$ws.Range("K131").Value = 71434152
$ws.Range("M131").Value = -71429112

$ws.Range("H137").Value = 1211498.4
$ws.Range("I137").Value = 3757.1538
$ws.Range("K137").Value = 11271.4614
$ws.Range("M137").Value = -8721.4614

$ws.Range("H138").Value = 3092996
$ws.Range("J138").Value = 4835767.5
$ws.Range("L138").Value = 14507302.5
$ws.Range("N138").Value = -14517582.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1048218.44
$ws.Range("I61").Value = 6862.2
$ws.Range("K61").Value = 6862.2
$ws.Range("M61").Value = -6650.2

$ws.Range("H63").Value = 5011.2666
$ws.Range("J63").Value = 6597.5557
$ws.Range("L63").Value = 6597.5557
$ws.Range("N63").Value = -7969.5557

$ws.Range("H66").Value = 5011.2666
$ws.Range("J66").Value = 6597.5557
$ws.Range("L66").Value = 32987.7785
$ws.Range("N66").Value = -39851.7785

$ws.Range("H74").Value = 2041.2
$ws.Range("I74").Value = 1782.0952
$ws.Range("J74").Value = 3401.5
$ws.Range("K74").Value = 1782.0952
$ws.Range("L74").Value = 3401.5
$ws.Range("M74").Value = -908.0952
$ws.Range("N74").Value = -5149.5

$ws.Range("H77").Value = 2041.2
$ws.Range("I77").Value = 1782.0952
$ws.Range("J77").Value = 3401.5
$ws.Range("K77").Value = 8910.476000000001
$ws.Range("L77").Value = 17007.5
$ws.Range("M77").Value = -4542.476000000001
$ws.Range("N77").Value = -25743.5

$ws.Range("H102").Value = 3867.25
$ws.Range("I102").Value = 1823
$ws.Range("K102").Value = 1823
$ws.Range("M102").Value = -201

$ws.Range("H136").Value = 1048218.44
$ws.Range("I136").Value = 6862.2
$ws.Range("K136").Value = 20586.6
$ws.Range("M136").Value = -18036.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 38466464
$ws.Range("I86").Value = 6285.7896
$ws.Range("K86").Value = 6285.7896
$ws.Range("M86").Value = -5162.7896

$ws.Range("H89").Value = 38466464
$ws.Range("I89").Value = 6285.7896
$ws.Range("K89").Value = 31428.948
$ws.Range("M89").Value = -25812.948

$ws.Range("H94").Value = 1899.2354
$ws.Range("I94").Value = 1837.8
$ws.Range("K94").Value = 1837.8
$ws.Range("M94").Value = -1386.8

$ws.Range("H99").Value = 2442.8333
$ws.Range("I99").Value = 1581.4
$ws.Range("K99").Value = 1581.4
$ws.Range("M99").Value = -83.40000000000009

$ws.Range("H105").Value = 4926
$ws.Range("I105").Value = 5581.913
$ws.Range("J105").Value = 2770.8572
$ws.Range("K105").Value = 5581.913
$ws.Range("L105").Value = 2770.8572
$ws.Range("M105").Value = -3834.913
$ws.Range("N105").Value = -6264.8572

$ws.Range("H134").Value = 1044495.06
$ws.Range("I134").Value = 2816.5293
$ws.Range("K134").Value = 8449.5879
$ws.Range("M134").Value = -5914.5879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 14304.8
$ws.Range("I103").Value = 14304.8
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 14304.8
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -13132.8
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 30305128
$ws.Range("J122").Value = 111112280
$ws.Range("L122").Value = 333336840
$ws.Range("N122").Value = -333341740

$ws.Range("H132").Value = 2282.7778
$ws.Range("I132").Value = 2413.5715
$ws.Range("K132").Value = 7240.7145
$ws.Range("M132").Value = -4710.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 5285.7144
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 6000
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 18000
$ws.Range("M57").Value = -2441
$ws.Range("N57").Value = -19118

$ws.Range("H131").Value = 3763935
$ws.Range("I131").Value = 200939.8
$ws.Range("J131").Value = 5134317.5
$ws.Range("K131").Value = 602819.3999999999
$ws.Range("L131").Value = 15402952.5
$ws.Range("M131").Value = -597779.3999999999
$ws.Range("N131").Value = -15413032.5

$ws.Range("H137").Value = 14615.6
$ws.Range("J137").Value = 19808.285
$ws.Range("L137").Value = 59424.855
$ws.Range("N137").Value = -69624.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H80").Value = 62502880
$ws.Range("I80").Value = 105264590
$ws.Range("J80").Value = 4988.6924
$ws.Range("K80").Value = 105264590
$ws.Range("L80").Value = 4988.6924
$ws.Range("M80").Value = -105263592
$ws.Range("N80").Value = -6984.6924

$ws.Range("H83").Value = 62502880
$ws.Range("I83").Value = 105264590
$ws.Range("J83").Value = 4988.6924
$ws.Range("K83").Value = 526322950
$ws.Range("L83").Value = 24943.462
$ws.Range("M83").Value = -526317958
$ws.Range("N83").Value = -34927.462

$ws.Range("H97").Value = 2213.0908
$ws.Range("I97").Value = 2548
$ws.Range("K97").Value = 2548
$ws.Range("M97").Value = -2052

$ws.Range("H102").Value = 1346.9354
$ws.Range("I102").Value = 1119.8214
$ws.Range("K102").Value = 1119.8214
$ws.Range("M102").Value = 502.1786

$ws.Range("H107").Value = 1324.25
$ws.Range("I107").Value = 564.3570999999999
$ws.Range("J107").Value = 1915.2778
$ws.Range("K107").Value = 564.3570999999999
$ws.Range("L107").Value = 1915.2778
$ws.Range("M107").Value = 1355.6429
$ws.Range("N107").Value = -5755.2778

$ws.Range("H122").Value = 1965.1923
$ws.Range("I122").Value = 1594.4
$ws.Range("K122").Value = 4783.200000000001
$ws.Range("M122").Value = -2333.200000000001

$ws.Range("H126").Value = 3010.7778
$ws.Range("I126").Value = 2499.7144
$ws.Range("J126").Value = 4799.5
$ws.Range("K126").Value = 7499.1432
$ws.Range("L126").Value = 14398.5
$ws.Range("M126").Value = -5029.1432
$ws.Range("N126").Value = -19338.5

$ws.Range("H132").Value = 43480996
$ws.Range("I132").Value = 58824540
$ws.Range("J132").Value = 7610.3335
$ws.Range("K132").Value = 176473620
$ws.Range("L132").Value = 22831.0005
$ws.Range("M132").Value = -176471090
$ws.Range("N132").Value = -27891.0005

$ws.Range("H135").Value = 169999.5
$ws.Range("J135").Value = 169999.5
$ws.Range("L135").Value = 169999.5
$ws.Range("N135").Value = -180139.5

$ws.Range("H136").Value = 19525.035
$ws.Range("J136").Value = 19525.035
$ws.Range("L136").Value = 58575.105
$ws.Range("N136").Value = -63675.105

$ws.Range("H138").Value = 87000
$ws.Range("J138").Value = 87000
$ws.Range("L138").Value = 87000
$ws.Range("N138").Value = -97280

$ws.Range("H139").Value = 107316.164
$ws.Range("J139").Value = 107316.164
$ws.Range("L139").Value = 107316.164
$ws.Range("N139").Value = -117596.164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 55558508
$ws.Range("I82").Value = 1109
$ws.Range("J82").Value = 125005256
$ws.Range("K82").Value = 1109
$ws.Range("L82").Value = 125005256
$ws.Range("M82").Value = -748
$ws.Range("N82").Value = -125005978

$ws.Range("H85").Value = 55558508
$ws.Range("I85").Value = 1109
$ws.Range("J85").Value = 125005256
$ws.Range("K85").Value = 1109
$ws.Range("L85").Value = 125005256
$ws.Range("M85").Value = 139
$ws.Range("N85").Value = -125007752

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 14499.25
$ws.Range("J55").Value = 14499.25
$ws.Range("L55").Value = 14499.25
$ws.Range("N55").Value = -15053.25

$ws.Range("H96").Value = 3641.9443
$ws.Range("I96").Value = 2505.5
$ws.Range("K96").Value = 2505.5
$ws.Range("M96").Value = -1132.5

$ws.Range("H107").Value = 22727650
$ws.Range("I107").Value = 334.2353
$ws.Range("J107").Value = 100000520
$ws.Range("K107").Value = 1002.7059
$ws.Range("L107").Value = 300001560
$ws.Range("M107").Value = 917.2941000000001
$ws.Range("N107").Value = -300005400

$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360
